$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 6
$ws.Range("F5").Value = 437
$ws.Range("F6").Value = 1278
$ws.Range("F8").Value = 7581
$ws.Range("F12").Value = 8217
$ws.Range("F15").Value = 60
$ws.Range("F16").Value = 5622
$ws.Range("F17").Value = 5622
$ws.Range("F19").Value = 2574
$ws.Range("F20").Value = 1109
$ws.Range("F21").Value = 4583
$ws.Range("F26").Value = 498
$ws.Range("F27").Value = 3099
$ws.Range("F28").Value = 3099
$ws.Range("F30").Value = 10
$ws.Range("F31").Value = 2855
$ws.Range("F32").Value = 2855
$ws.Range("F34").Value = 326
$ws.Range("F35").Value = 122
$ws.Range("F36").Value = 284
$ws.Range("F38").Value = 639
$ws.Range("F40").Value = 871
$ws.Range("F41").Value = 1634
$ws.Range("F44").Value = 5
$ws.Range("F45").Value = 2647
$ws.Range("F47").Value = 2272
$ws.Range("B50").Value = "'2024-05-25"
$ws.Range("E50").Value = "2024.05.25 09:00-05.26 17:00"
$ws.Range("F50").Value = 470
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 39
$ws.Range("F7").Value = 32
$ws.Range("F8").Value = 106
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 1311
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1311
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 1278
$ws.Range("F7").Value = 7581
$ws.Range("F11").Value = 8217
$ws.Range("F13").Value = 60
$ws.Range("F14").Value = 5622
$ws.Range("F15").Value = 5622
$ws.Range("F17").Value = 2574
$ws.Range("F18").Value = 1109
$ws.Range("F19").Value = 4583
$ws.Range("F20").Value = 0
$ws.Range("F25").Value = 498
$ws.Range("F26").Value = 3099
$ws.Range("F27").Value = 3100
$ws.Range("F29").Value = 10
$ws.Range("F30").Value = 2855
$ws.Range("F31").Value = 2855
$ws.Range("F32").Value = 326
$ws.Range("F33").Value = 122
$ws.Range("F34").Value = 284
$ws.Range("F35").Value = 39
$ws.Range("F37").Value = 639
$ws.Range("F40").Value = 871
$ws.Range("F42").Value = 1634
$ws.Range("F45").Value = 5
$ws.Range("F46").Value = 2647
$ws.Range("F47").Value = 32
$ws.Range("F48").Value = 2272
$ws.Range("B51").Value = "'2024-05-25"
$ws.Range("E51").Value = "2024.05.25 09:00-05.26 17:00"
$ws.Range("F51").Value = 470
$ws.Range("F52").Value = 106
